$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "67.445.07"
Set-TextValue "E2" "  -3.33%  "
Set-TextValue "D3" "3.707.74"
Set-TextValue "E3" "  -3.91%  "
Set-TextValue "D4" "1.00"
Set-TextValue "E4" "  +0.08%  "
Set-TextValue "D5" "596.50"
Set-TextValue "E5" "  -2.09%  "
Set-TextValue "D6" "165.79"
Set-TextValue "E6" "  -5.04%  "
Set-TextValue "D7" "3.707.90"
Set-TextValue "E7" "  -4.04%  "
Set-TextValue "E8" "  +0.07%  "
Set-TextValue "D9" "0.531"
Set-TextValue "E9" "  +0.81%  "
Set-TextValue "D10" "0.162"
Set-TextValue "E10" "  -2.60%  "
Set-TextValue "D11" "6.19"
Set-TextValue "E11" "  -4.55%  "
Set-TextValue "D12" "0.462"
Set-TextValue "E12" "  -3.77%  "
Set-TextValue "D13" "37.69"
Set-TextValue "E13" "  -5.56%  "
Set-TextValue "E14" "  -4.80%  "
Set-TextValue "D15" "4.328.32"
Set-TextValue "E15" "  -3.67%  "
Set-TextValue "D16" "3.710.75"
Set-TextValue "E16" "  -3.53%  "
Set-TextValue "D17" "67.522.92"
Set-TextValue "E17" "  -3.29%  "
Set-TextValue "D18" "17.56"
Set-TextValue "E18" "  +5.72%  "
Set-TextValue "E19" "  -3.48%  "
Set-TextValue "E20" "  -3.00%  "
Set-TextValue "D21" "492.02"
Set-TextValue "E21" "  -2.66%  "
Set-TextValue "D22" "9.31"
Set-TextValue "E22" "  -1.90%  "
Set-TextValue "D23" "0.724"
Set-TextValue "E23" "  -2.05%  "
Set-TextValue "D24" "85.76"
Set-TextValue "E24" "  -0.19%  "
Set-TextValue "D25" "2.30"
Set-TextValue "E25" "  -5.94%  "
Set-TextValue "E26" "  -3.02%  "
Set-TextValue "E27" "  -3.18%  "
Set-TextValue "D28" "10.11"
Set-TextValue "E28" "  -3.15%  "
Set-TextValue "E29" "  +0.14%  "
Set-TextValue "D30" "2.94"
Set-TextValue "E30" "  -1.52%  "
Set-TextValue "E31" "  -6.57%  "
Set-TextValue "D32" "31.53"
Set-TextValue "E32" "  -2.87%  "
Set-TextValue "D33" "7.63"
Set-TextValue "E33" "  -3.77%  "
Set-TextValue "D34" "3.844.46"
Set-TextValue "E34" "  -3.76%  "
Set-TextValue "E35" "  -4.32%  "
Set-TextValue "D36" "3.650.72"
Set-TextValue "E36" "  -3.67%  "
Set-TextValue "E37" "  +0.03%  "
Set-TextValue "D38" "0.996"
Set-TextValue "E38" "  -4.85%  "
Set-TextValue "E39" "  -5.61%  "
Set-TextValue "E40" "  -6.75%  "
Set-TextValue "D41" "0.322"
Set-TextValue "E41" "  -3.61%  "
Set-TextValue "D42" "433.16"
Set-TextValue "E42" "  -10.07%  "
Set-TextValue "D43" "48.59"
Set-TextValue "E43" "  -2.49%  "
Set-TextValue "E44" "  -5.25%  "
Set-TextValue "D45" "2.79"
Set-TextValue "E45" "  -6.61%  "
Set-TextValue "E46" "  -1.25%  "
Set-TextValue "D47" "40.74"
Set-TextValue "E47" "  -5.64%  "
Set-TextValue "E48" "  +0.02%  "
Set-TextValue "D49" "143.03"
Set-TextValue "E49" "  +2.12%  "
Set-TextValue "D50" "2.755.03"
Set-TextValue "E50" "  -5.53%  "
Set-TextValue "E51" "  -3.32%  "
